$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (row 1): A..T ---
$headers = @(
    "Places",
    "location",
    "Temperature1",
    "Temperature2",
    "Temperature3",
    "Temperature4",
    "Temperature5",
    "Temperature6",
    "Temperature7",
    "Temperature8",
    "Temperature9",
    "Temperature10",
    "Temperature11",
    "Temperature12",
    "Temperature13",
    "Temperature14",
    "Temperature15",
    "Temperature16",
    "Temperature17",
    "Temperature18"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
# keep the existing yellow header fill consistent across the newly added columns
$ws.Range("A1:T1").Interior.Color = 65535

# --- Clear out old data rows/columns beyond what the new layout needs ---
$ws.Range("A2:T100").Clear()

# --- New data rows (row 2..5): Places / location (B & C) / Temperature16 (R & S) ---
$places = @("dusseldorf", "Nice", "Marseille", "Monte Carlo")
$locations = @(
    "16.04.`$[[ChromeDriver: chrome on XP (70d13a5d04fcbe4353588bb598b793d3)] -> id: xPat]",
    "13.58.`$[[ChromeDriver: chrome on XP (6467ec0d74fea02c0a3ac09ad57bcefb)] -> id: xPat]",
    "20.36.`$[[ChromeDriver: chrome on XP (d2699a01bbe2d02d9fdb2c4a0859a00f)] -> id: xPat]",
    "21.96.`$[[ChromeDriver: chrome on XP (52ab270f72e9e43db6022c5f290b4b9d)] -> id: xPat]"
)
$temp16 = @("13.3.", "13.28.", "19.12.", "21.81.")

for ($i = 0; $i -lt $places.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $places[$i]
    $ws.Cells.Item($row, 2).Value = $locations[$i]
    $ws.Cells.Item($row, 3).Value = $locations[$i]
    $ws.Cells.Item($row, 18).Value = $temp16[$i]
    $ws.Cells.Item($row, 19).Value = $temp16[$i]
}

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 83.03125
$ws.Columns.Item(3).ColumnWidth = 83.03125
for ($c = 4; $c -le 11; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 13.85546875
}
for ($c = 12; $c -le 19; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 15.0
}

# --- Selection matches the target (B1 active) ---
$ws.Range("B1").Select() | Out-Null
